# Risk Management Report - add Week4.1 risk entry (row 14)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Week4.1"
$ws.Range("B14").Value = "Game doesn't end when player has no card"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = "only limit to 13 rounds"
$ws.Range("G14").Value = "Yes"

# Move the active selection to A15, as in the saved workbook after data entry
$ws.Range("A15").Select() | Out-Null
